# Jay Tucker Mendoza Q0488 - Training Dashboard updates
# "adding averages and more checks": refresh the "last update" check date and
# the recomputed "period to expire" counters, and restyle the header/title
# text to bold white (on the existing dark-blue banner) instead of the old
# bold-black / bold-14pt look.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Re-run the periodic checks: new check date + refreshed day counts ---
$ws.Range("H3").Value = 113
$ws.Range("I3").Value = "'16-Sep-2025"

$ws.Range("H4").Value = 680
$ws.Range("I4").Value = "'16-Sep-2025"

$ws.Range("H5").Value = 679
$ws.Range("I5").Value = "'16-Sep-2025"

$ws.Range("H6").Value = 680
$ws.Range("I6").Value = "'16-Sep-2025"

$ws.Range("H7").Value = 679
$ws.Range("I7").Value = "'16-Sep-2025"

# --- Restyle: bold white text for the title and the column header band ---
# The title/header styles are shared across both sheets in the workbook's
# style table, so both dashboards get the refreshed bold-white look.
foreach ($sheet in $wb.Worksheets) {
    $usedRange = $sheet.UsedRange
    $lastCol = $usedRange.Columns.Count
    $lastColLetter = [char](64 + $lastCol)

    $titleRange = $sheet.Range("A1")
    $titleRange.Font.Bold = $true
    $titleRange.Font.Size = 11
    $titleRange.Font.Color = 16777215

    $headerRange = $sheet.Range("A2:" + $lastColLetter + "2")
    $headerRange.Font.Bold = $true
    $headerRange.Font.Size = 11
    $headerRange.Font.Color = 16777215
}
